$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 341.16666
$ws.Range("I18").Value = 341.16666
$ws.Range("K18").Value = 341.16666
$ws.Range("M18").Value = -57.16665999999998

$ws.Range("H40").Value = 2358.1
$ws.Range("I40").Value = 2750.25
$ws.Range("J40").Value = 2096.6667
$ws.Range("K40").Value = 2750.25
$ws.Range("L40").Value = 2096.6667
$ws.Range("M40").Value = -2575.25
$ws.Range("N40").Value = -2446.6667

$ws.Range("H41").Value = 6944676
$ws.Range("J41").Value = 290
$ws.Range("L41").Value = 290
$ws.Range("N41").Value = -1170

$ws.Range("H64").Value = 836938.75
$ws.Range("I64").Value = 912569.25
$ws.Range("J64").Value = 5003
$ws.Range("K64").Value = 912569.25
$ws.Range("L64").Value = 5003
$ws.Range("M64").Value = -912321.25
$ws.Range("N64").Value = -5499

$ws.Range("H67").Value = 836938.75
$ws.Range("I67").Value = 912569.25
$ws.Range("J67").Value = 5003
$ws.Range("K67").Value = 912569.25
$ws.Range("L67").Value = 5003
$ws.Range("M67").Value = -911711.25
$ws.Range("N67").Value = -6719

$ws.Range("H98").Value = 622707.25
$ws.Range("I98").Value = 746622
$ws.Range("J98").Value = 3133.3333
$ws.Range("K98").Value = 746622
$ws.Range("L98").Value = 3133.3333
$ws.Range("M98").Value = -745124
$ws.Range("N98").Value = -6129.3333

$ws.Range("H106").Value = 22224922
$ws.Range("I106").Value = 22224922
$ws.Range("K106").Value = 22224922
$ws.Range("M106").Value = -22224291

$ws.Range("H111").Value = 971.4
$ws.Range("I111").Value = 839.25
$ws.Range("J111").Value = 1500
$ws.Range("K111").Value = 2517.75
$ws.Range("L111").Value = 4500
$ws.Range("M111").Value = 549.25
$ws.Range("N111").Value = -10634

$ws.Range("H122").Value = 622707.25
$ws.Range("I122").Value = 746622
$ws.Range("J122").Value = 3133.3333
$ws.Range("K122").Value = 2239866
$ws.Range("L122").Value = 9399.999899999999
$ws.Range("M122").Value = -2237416
$ws.Range("N122").Value = -14299.9999

$ws.Range("H123").Value = 99988
$ws.Range("J123").Value = 99988
$ws.Range("L123").Value = 99988
$ws.Range("N123").Value = -109788

$ws.Range("H138").Value = 6414874.5
$ws.Range("I138").Value = 2305103
$ws.Range("K138").Value = 6915309
$ws.Range("M138").Value = -6910169

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2533.2354
$ws.Range("I61").Value = 1564.5454
$ws.Range("J61").Value = 4309.1665
$ws.Range("K61").Value = 1564.5454
$ws.Range("L61").Value = 4309.1665
$ws.Range("M61").Value = -1352.5454
$ws.Range("N61").Value = -4733.1665

$ws.Range("H123").Value = 35214
$ws.Range("J123").Value = 35214
$ws.Range("L123").Value = 35214
$ws.Range("N123").Value = -45014

$ws.Range("H136").Value = 2533.2354
$ws.Range("I136").Value = 1564.5454
$ws.Range("J136").Value = 4309.1665
$ws.Range("K136").Value = 4693.6362
$ws.Range("L136").Value = 12927.4995
$ws.Range("M136").Value = -2143.6362
$ws.Range("N136").Value = -18027.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2917.054
$ws.Range("I134").Value = 1869.6333
$ws.Range("K134").Value = 5608.8999
$ws.Range("M134").Value = -3073.8999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1547.1
$ws.Range("I31").Value = 1108.2
$ws.Range("J31").Value = 3741.6
$ws.Range("K31").Value = 1108.2
$ws.Range("L31").Value = 3741.6
$ws.Range("M31").Value = -813.2
$ws.Range("N31").Value = -4331.6

$ws.Range("H34").Value = 1547.1
$ws.Range("I34").Value = 1108.2
$ws.Range("J34").Value = 3741.6
$ws.Range("K34").Value = 1108.2
$ws.Range("L34").Value = 3741.6
$ws.Range("M34").Value = -906.2
$ws.Range("N34").Value = -4145.6

$ws.Range("H132").Value = 2846.9473
$ws.Range("I132").Value = 1519.4546
$ws.Range("J132").Value = 4672.25
$ws.Range("K132").Value = 4558.3638
$ws.Range("L132").Value = 14016.75
$ws.Range("M132").Value = -2028.3638
$ws.Range("N132").Value = -19076.75

$ws.Range("H134").Value = 2967.92
$ws.Range("I134").Value = 1418.5
$ws.Range("J134").Value = 5722.4443
$ws.Range("K134").Value = 4255.5
$ws.Range("L134").Value = 17167.3329
$ws.Range("M134").Value = -1720.5
$ws.Range("N134").Value = -22237.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 51500
$ws.Range("I3").Value = 51500
$ws.Range("K3").Value = 154500
$ws.Range("M3").Value = -154388

$ws.Range("H113").Value = 50001532
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 50001532
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 150004596
$ws.Range("N113").Value = -150008936
$ws.Range("M113").ClearContents()

$ws.Range("H122").Value = 938.375
$ws.Range("J122").Value = 983.8333
$ws.Range("L122").Value = 8854.4997
$ws.Range("N122").Value = -13754.4997

$ws.Range("H131").Value = 2751.9167
$ws.Range("I131").Value = 320
$ws.Range("J131").Value = 2857.652
$ws.Range("K131").Value = 960
$ws.Range("L131").Value = 8572.956
$ws.Range("M131").Value = 4080
$ws.Range("N131").Value = -18652.956

$ws.Range("H132").Value = 1254.5
$ws.Range("I132").Value = 1112.4
$ws.Range("J132").Value = 1396.6
$ws.Range("K132").Value = 10011.6
$ws.Range("L132").Value = 12569.4
$ws.Range("M132").Value = -7481.6
$ws.Range("N132").Value = -17629.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 10170
$ws.Range("J123").Value = 10170
$ws.Range("L123").Value = 10170
$ws.Range("N123").Value = -15070

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2825.1738
$ws.Range("I100").Value = 2397
$ws.Range("J100").Value = 3012.5
$ws.Range("K100").Value = 2397
$ws.Range("L100").Value = 3012.5
$ws.Range("M100").Value = -1856
$ws.Range("N100").Value = -4094.5

$ws.Range("H123").Value = 36000
$ws.Range("J123").Value = 36000
$ws.Range("L123").Value = 36000
$ws.Range("N123").Value = -45800

$ws.Range("H136").Value = 4869.0713
$ws.Range("I136").Value = 1724.9375
$ws.Range("J136").Value = 9061.25
$ws.Range("K136").Value = 5174.8125
$ws.Range("L136").Value = 27183.75
$ws.Range("M136").Value = -2624.8125
$ws.Range("N136").Value = -32283.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 39210
$ws.Range("J123").Value = 39210
$ws.Range("L123").Value = 39210
$ws.Range("N123").Value = -49010

$ws.Range("H136").Value = 15922108
$ws.Range("I136").Value = 19667796
$ws.Range("K136").Value = 59003388
$ws.Range("M136").Value = -59000838
